$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Statut value (from Curia Vista affair status)
$statutMap = @{}
$statutMap[2] = 'Überwiesen an den Bundesrat / Transmis au Conseil fédéral'
$statutMap[3] = 'Erledigt / Liquidé'
$statutMap[4] = 'In Kommission des Ständerats / En commission du Conseil des Etats'
$statutMap[5] = 'Stellungnahme zum Vorstoss liegt vor / L’avis relatif à l’intervention est disponible'
$statutMap[6] = 'Stellungnahme zum Vorstoss liegt vor / L’avis relatif à l’intervention est disponible'
$statutMap[7] = 'Erledigt / Liquidé'
$statutMap[8] = 'Erledigt / Liquidé'
$statutMap[9] = 'Zugewiesen an die behandelnde Kommission / Attribué à la commission compétente'
$statutMap[10] = 'Stellungnahme zum Vorstoss liegt vor / L’avis relatif à l’intervention est disponible'
$statutMap[11] = 'Erledigt / Liquidé'
$statutMap[12] = 'Stellungnahme zum Vorstoss liegt vor / L’avis relatif à l’intervention est disponible'
$statutMap[13] = 'Stellungnahme zum Vorstoss liegt vor / L’avis relatif à l’intervention est disponible'
$statutMap[14] = 'Eingereicht / Déposé'
$statutMap[15] = 'Eingereicht / Déposé'
$statutMap[16] = 'Eingereicht / Déposé'
$statutMap[17] = 'Eingereicht / Déposé'
$statutMap[18] = 'Erledigt / Liquidé'
$statutMap[19] = 'Erledigt / Liquidé'
$statutMap[20] = 'Erledigt / Liquidé'
$statutMap[21] = 'Erledigt / Liquidé'
$statutMap[22] = 'Erledigt / Liquidé'
$statutMap[23] = 'Erledigt / Liquidé'
$statutMap[24] = 'Erledigt / Liquidé'
$statutMap[25] = 'Erledigt / Liquidé'

# Grow the table to a 6th column first, so the ListObject/table
# definition tracks the extra column from the start.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F25"))

# Shift the existing Lien_FR column (E) data into the new column F
$ws.Range("E1:E25").Copy($ws.Range("F1:F25"))

# Write the new Statut values into column E
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = $statutMap[$r]
}

# Relabel the Type column: every exported row here is a content update
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 1).Value = "Mise à jour contenu"
}

# Refresh the table header names through the ListObject header range
# (keeps the table parts cached column names in sync with the cells).
$lo.HeaderRowRange.Cells.Item(1, 1).Value = "Type_Changement"
$lo.HeaderRowRange.Cells.Item(1, 2).Value = "Numéro"
$lo.HeaderRowRange.Cells.Item(1, 3).Value = "Auteur"
$lo.HeaderRowRange.Cells.Item(1, 4).Value = "Mention"
$lo.HeaderRowRange.Cells.Item(1, 5).Value = "Statut"
$lo.HeaderRowRange.Cells.Item(1, 6).Value = "Lien_FR"
